# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 172 (pushing existing rows 172-198 down
# to 173-199) and populate it with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 172:198 down to 173:199, creating a blank row 172.
$ws.Rows.Item(172).Insert()

# Populate the new row 172 with the new weekly record.
$ws.Range("A172").Value = 11
$ws.Range("B172").Value = "Vega Monumental Concepción"
$ws.Range("C172").Value = "Bíobío"
$ws.Range("D172").Value = 45218
$ws.Range("E172").Value = 8
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100108
$ws.Range("H172").Value = "Tropicales y subtropicales"
$ws.Range("I172").Value = 100108002
$ws.Range("J172").Value = "Mango"
$ws.Range("K172").Value = "Sin especificar"
$ws.Range("L172").Value = "Primera"
$ws.Range("M172").Value = 100
$ws.Range("N172").Value = 11000
$ws.Range("O172").Value = 11000
$ws.Range("P172").Value = 11000
$ws.Range("Q172").Value = "`$/bandeja 4 kilos"
$ws.Range("R172").Value = "Brasil"
$ws.Range("S172").Value = 2750
$ws.Range("T172").Value = 4

# Make sure the style on the date cell matches the other date cells in the
# column (carry the same format used by D173:D199).
$ws.Range("D172").NumberFormat = $ws.Range("D173").NumberFormat
